# Apply the updated "Fitness" (column C) values for run_13 log rows, as
# captured by the source diff. Column A (Run) and B (Generation) are left
# untouched; only the Fitness values in column C change, in contiguous
# "plateau" blocks that match consecutive rows sharing the same new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C39").Value = 8342
$ws.Range("C40:C50").Value = 8322
$ws.Range("C51:C57").Value = 8283
$ws.Range("C58:C83").Value = 8106
$ws.Range("C84:C86").Value = 7345
$ws.Range("C87:C97").Value = 7295
